$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ06422413",
    "summ06505780",
    "summ06580670",
    "summ06663159",
    "summ06752403",
    "summ06823239",
    "summ06924734",
    "summ07008671",
    "summ07090858",
    "summ07173239",
    "summ07256533",
    "summ07355120",
    "summ07424589",
    "summ07516593",
    "summ07591366",
    "summ07674601",
    "summ07761034",
    "summ07841436",
    "summ07924285",
    "summ08031244",
    "summ08109646",
    "summ08190881",
    "summ08255487",
    "summ08340648",
    "summ08423089",
    "summ08507619",
    "summ08605725",
    "summ08690738",
    "summ08772445",
    "summ08848888",
    "summ08923077",
    "summ09007206",
    "summ09090660",
    "summ09176582",
    "summ09255743",
    "summ09345515",
    "summ09429396",
    "summ09514496",
    "summ09599198",
    "summ09682382",
    "summ09757056",
    "summ09841314",
    "summ09923207",
    "summ10005195",
    "summ10082301",
    "summ10165618",
    "summ10240372",
    "summ10323251",
    "summ10405067",
    "summ10482382",
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
